# Weekly update: a new "Coco" price record for Vega Modelo de Temuco is
# inserted at row 29, pushing all existing records down by one row
# (old row 29 becomes row 30, ..., old row 116 becomes row 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 29 (shifts rows 29..116 down to 30..117)
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new weekly record
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 45148
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100108
$ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(29, 9).Value = 100108007
$ws.Cells.Item(29, 10).Value = "Coco"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 100
$ws.Cells.Item(29, 14).Value = 36000
$ws.Cells.Item(29, 15).Value = 36000
$ws.Cells.Item(29, 16).Value = 36000
$ws.Cells.Item(29, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(29, 18).Value = "Perú"
$ws.Cells.Item(29, 19).Value = 1800
$ws.Cells.Item(29, 20).Value = 20
